# Apply the "gaps_benchmark/all1.xlsx" refresh:
#  - insert a new "wait-and-see" results column between the existing
#    FA-MSP-R and RH2SSP-R columns (both in the raw-value block A:E and
#    in the gap-percentage block G:J)
#  - refresh all of the underlying values with the re-run numbers
#  - re-enter the gap formulas so they pick up the new column and use a
#    one-decimal percentage format
#  - tidy up the selection to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Headers: A/B stay put, a new "wait-and-see" column lands at C (and
#    H in the mirrored gap block); the old RH2SSP-R / static2SSP-R
#    headers shift one slot to the right.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "wait-and-see"
$ws.Range("D1").Value = "RH2SSP-R"
$ws.Range("E1").Value = "static2SSP-R"

$ws.Range("H1").Value = "wait-and-see"
$ws.Range("I1").Value = "RH2SSP-R"
$ws.Range("J1").Value = "static2SSP-R"

# ---------------------------------------------------------------------
# 2. Raw values (columns A, B, C, D, E) and the probability column (F)
#    Column E used to be an empty, centre-aligned placeholder (E2:E11);
#    now it carries real numbers with the plain default style, so drop
#    the old formatting before writing into it.
# ---------------------------------------------------------------------
$ws.Range("E2:E11").ClearFormats()

$ws.Range("A2").Value = 637.86344150000002
$ws.Range("B2").Value = 651.73575530000005
$ws.Range("C2").Value = 651.43995580000001
$ws.Range("D2").Value = 4094.2388740000001
$ws.Range("E2").Value = 4838.6550889999999
$ws.Range("F2").Value = 0.05

$ws.Range("A3").Value = 655.4050416
$ws.Range("B3").Value = 764.19803730000001
$ws.Range("C3").Value = 766.70473140000001
$ws.Range("D3").Value = 4128.7297189999999
$ws.Range("E3").Value = 4849.1916929999998
$ws.Range("F3").Value = 0.1

$ws.Range("A4").Value = 665.97388920000003
$ws.Range("B4").Value = 875.14518090000001
$ws.Range("C4").Value = 881.96950700000002
$ws.Range("D4").Value = 4162.3372390000004
$ws.Range("E4").Value = 4860.0398519999999
$ws.Range("F4").Value = 0.15

$ws.Range("A5").Value = 676.40887650000002
$ws.Range("B5").Value = 980.7120668
$ws.Range("C5").Value = 997.23428249999995
$ws.Range("D5").Value = 4186.6777339999999
$ws.Range("E5").Value = 4873.8786140000002
$ws.Range("F5").Value = 0.2

$ws.Range("A6").Value = 686.84386380000001
$ws.Range("B6").Value = 1069.965144
$ws.Range("C6").Value = 1112.4990580000001
$ws.Range("D6").Value = 4210.8534730000001
$ws.Range("E6").Value = 4881.0714170000001
$ws.Range("F6").Value = 0.25

$ws.Range("A7").Value = 697.27885119999996
$ws.Range("B7").Value = 1153.539538
$ws.Range("C7").Value = 1227.7638340000001
$ws.Range("D7").Value = 4248.4316790000003
$ws.Range("E7").Value = 4891.4158269999998
$ws.Range("F7").Value = 0.3

$ws.Range("A8").Value = 707.71383849999995
$ws.Range("B8").Value = 1231.56807
$ws.Range("C8").Value = 1343.028609
$ws.Range("D8").Value = 4278.3880680000002
$ws.Range("E8").Value = 4901.8197929999997
$ws.Range("F8").Value = 0.35

$ws.Range("A9").Value = 718.14882580000005
$ws.Range("B9").Value = 1306.806243
$ws.Range("C9").Value = 1458.2933849999999
$ws.Range("D9").Value = 4310.1987669999999
$ws.Range("E9").Value = 4912.3304969999999
$ws.Range("F9").Value = 0.4

$ws.Range("A10").Value = 728.58381320000001
$ws.Range("B10").Value = 1379.487171
$ws.Range("C10").Value = 1573.55816
$ws.Range("D10").Value = 4316.1801260000002
$ws.Range("E10").Value = 4922.7347719999998
$ws.Range("F10").Value = 0.45

$ws.Range("A11").Value = 739.0188005
$ws.Range("B11").Value = 1446.4878080000001
$ws.Range("C11").Value = 1688.822936
$ws.Range("D11").Value = 4339.8287280000004
$ws.Range("E11").Value = 4933.1546669999998
$ws.Range("F11").Value = 0.5

# ---------------------------------------------------------------------
# 3. Gap-vs-baseline formulas (G:J), one-decimal percent format
# ---------------------------------------------------------------------
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("G$r").Formula = "=(B$r-`$A$r)/`$A$r"
    $ws.Range("H$r").Formula = "=(C$r-`$A$r)/`$A$r"
    $ws.Range("I$r").Formula = "=(D$r-`$A$r)/`$A$r"
    $ws.Range("J$r").Formula = "=(E$r-`$A$r)/`$A$r"
    $ws.Range("G$r`:J$r").NumberFormat = "0.0%"
}

# ---------------------------------------------------------------------
# 4. Tidy up the now-unused trailing placeholder cell in the old
#    9th-column formatting strip for rows 16 downward (the refreshed
#    export only keeps the F:H placeholder strip there).
# ---------------------------------------------------------------------
for ($r = 16; $r -le 52; $r++) {
    $ws.Range("I$r").ClearContents()
    $ws.Range("I$r").ClearFormats()
}

# ---------------------------------------------------------------------
# 5. Selection, matching the saved state in the workbook
# ---------------------------------------------------------------------
$ws.Range("F18").Select()
